# ----------------------------------------------------------------------
# 1) Append the new 'Knärot' section (heading, body paragraphs, reference
#    list) right after the last paragraph of the document, which is the
#    'BILAGA 1 - Fridlysta arter' title paragraph.
# ----------------------------------------------------------------------
$d = $word.ActiveDocument

# --- new paragraph 1 of 13 (style: Heading1) ---
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Style = 'Heading1'
$nr = $newPara.Range
$nr.InsertAfter('Knärot – ekologi samt krav på livsmiljön')

# --- new paragraph 2 of 13 (style: Normal) ---
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Style = 'Normal'
$nr = $newPara.Range
$nr.InsertAfter('Knärot är fridlyst enligt 8 och 15 §§ artskyddsförordningen och klassad som sårbar (VU) enligt rödlistan 2020. Knärot är beroende av hög och jämn luftfuktighet i gamla, ostörda skogsmiljöer och är känslig för snabba förändringar av ljus-/vindförhållanden eller uttorkning. På grund av ett alltför intensivt skogsbruk har den minskat med 40 (25-50) % under de senaste 60 åren och i framtiden bedöms minskningstakten uppgå till 30 (20-40) %. Till följd av att arten har en dokumenterat högre minskningstakt iförhållande till sin generationstid än vad som tidigare varit känt (data från Riksskogstaxeringen) höjdes den till hotkategori sårbar (VU) i rödlistan 2020 (Artdatabanken, 2021).')

# --- new paragraph 3 of 13 (style: Normal) ---
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Style = 'Normal'
$nr = $newPara.Range
$nr.InsertAfter('Samuel Johnsons doktorsavhandling ')
$s2_1 = $nr.End
$nr.InsertAfter('“Retention Forestry as a Conservation Measure for Boreal Forest Ground Vegetation“')
$e2_1 = $nr.End
$nr.InsertAfter(' (SLU, Uppsala 2014) visar att det krävs väl tilltagna skyddszoner för att knärotens växtplatser inte ska ta skada av skogsbruksåtgärder i intilliggande områden: ')
$s2_3 = $nr.End
$nr.InsertAfter('“Study III shows that retention patches smaller than 0.5 ha do not lifeboat the sensitive forest herb G. repens, a species that depend on stable microclimatic conditions typical for intact forest stands.” ')
$e2_3 = $nr.End
$nr.InsertAfter('Vidare ')
$s2_5 = $nr.End
$nr.InsertAfter('“More sensitive forest species are not lifeboated in retention patches ranging from 0.05 to 0.5 ha (Papers II & III).”')
$e2_5 = $nr.End
$ir = $d.Range($s2_1 - 1, $e2_1 - 1)
$ir.Font.Italic = $true
$ir = $d.Range($s2_3 - 1, $e2_3 - 1)
$ir.Font.Italic = $true
$ir = $d.Range($s2_5 - 1, $e2_5 - 1)
$ir.Font.Italic = $true

# --- new paragraph 4 of 13 (style: Normal) ---
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Style = 'Normal'
$nr = $newPara.Range
$nr.InsertAfter('Johnsons (2014) rekommendation på minst 50 meters breda skyddszoner runt knärotens växtplatser motsvarar en areal på 0,78 hektar, vilket ligger i linje med andra studier som gjorts på känsliga skogsarter: ')
$s3_1 = $nr.End
$nr.InsertAfter('“In study III I also show that translocated specimens of G. repens survives well in mature forests at least 50 m from the nearest edge to an open area. Moreover, measures of temperature and humidity show that such distances from an open area is far enough to offer a microclimate that is more stable compared to what present in retention patches of around 0.1 ha. This means that the very centre of a circular patch with radius 50 m (equals a size of 0.78 ha) should offer conditions similar to interior forest and would perhaps be a suitable habitat for G. repens and similar species. Previous studies from both North America and Sweden have also concluded that patches between 0.5 and one ha are sufficient for preserving interior forest vegetation as well as sensitive lichens and bryophytes (de Graaf & Roberts 2009; Halpern et al. 2012; Rudolphi et al. 2014).”')
$e3_1 = $nr.End
$ir = $d.Range($s3_1 - 1, $e3_1 - 1)
$ir.Font.Italic = $true

# --- new paragraph 5 of 13 (style: Normal) ---
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Style = 'Normal'
$nr = $newPara.Range
$nr.InsertAfter('En nyligen publicerad vetenskaplig uppsats av Koelmeijer m.fl. (2022) inkluderar orkidén knärots skyddsbehov. I uppsatsen berörs problemet med uttorkning för växter, bl.a. för knärot, ett problem som blivit accentuerat på grund av den pågående klimatförändringen och torra somrar, t.ex. den exceptionellt torra sommaren 2018. I uppsatsen undersöks områden med tre olika avstånd från kalhyggeskant med avseende på skydd bl.a. för knärot. Det första området har avstånd upp till 20 m från hyggeskant (Strong edge effect), det andra 20 – 40 m från hyggeskant (Weak edge effect) och det tredje avser större avstånd från hyggeskant, där kanteffekten anses vara försumbar (Interior). Ett resultat var att man fann stor eller mycket stor uttorkningseffekt på känsliga och rödlistade skogsarter vid de kortare avstånden till hyggeskant, medan effekt av uttorkning inte konstaterades på större avstånd (Interior). För orkidén knärot fann man en rik förekomst (upp till 0,06 dm2/m2) på stort avstånd från hyggeskant (Interior), medan förekomsten var liten eller närmast försumbar i de områden som klassificerades som Weak edge effect respektive Strong edge effect. Arbetet påpekar att de allt oftare förekommande torra somrarna ger ytterligare skäl att utöka skyddsavståndet från hyggen till den fuktkrävande arten knärot (Koelmeijer m.fl., 2022).')

# --- new paragraph 6 of 13 (style: Normal) ---
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Style = 'Normal'
$nr = $newPara.Range
$nr.InsertAfter('Även Skogsstyrelsens egen vägledning för hänsyn till knärot ligger i linje med ovanstående forskningsstudier. Av vägledningen framgår det att för med hög sannolikhet kunna bevara befintliga förekomster krävs relativt stora avsättningar av uppvuxen skog med slutet och relativt tätt kronskikt. Som riktlinje kan krävas ett avstånd på 50 meter in från brynet för att vidmakthålla ett fungerande mikroklimat. Detta innebär att fristående hänsynsytor för många arter (kärlväxter, lavar och mossor) kan behöva ha en area överstigande 0,8 hektar (cirkelyta med radien 50 meter = 0,78 hektar) för att bibehålla lokalklimatet. Även ganska små förändringar i form av förändrade ljus- och fuktighetsförhållanden, till exempel till följd av gallring, kan leda till att arten försvinner till följd av konkurrens med mera ljuskrävande och snabbväxande arter (Skogsstyrelsen, 2022).')

# --- new paragraph 7 of 13 (style: Heading2) ---
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Style = 'Heading2'
$nr = $newPara.Range
$nr.InsertAfter('Referenser - knärot')

# --- new paragraph 8 of 13 (style: Normal) ---
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Style = 'Normal'
$nr = $newPara.Range
$nr.InsertAfter('de Graaf M & Roberts M.R., 2009. ')
$s7_1 = $nr.End
$nr.InsertAfter('Short-term response of the herbaceous layer within leave patches after harvest. ')
$e7_1 = $nr.End
$nr.InsertAfter('Forest Ecology and Management 257, 1014-1025')
$ir = $d.Range($s7_1 - 1, $e7_1 - 1)
$ir.Font.Italic = $true

# --- new paragraph 9 of 13 (style: Normal) ---
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Style = 'Normal'
$nr = $newPara.Range
$nr.InsertAfter('Halpern, C. B., Halaj, J., Evans, S. A., & Dovciak, M., 2012. ')
$s8_1 = $nr.End
$nr.InsertAfter('Level and pattern of overstory retention interact to shape long-term responses of understories to timber harvest. ')
$e8_1 = $nr.End
$nr.InsertAfter('Ecological Applications, 22, 2049-2064 ')
$ir = $d.Range($s8_1 - 1, $e8_1 - 1)
$ir.Font.Italic = $true

# --- new paragraph 10 of 13 (style: Normal) ---
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Style = 'Normal'
$nr = $newPara.Range
$nr.InsertAfter('Koelmeijer, I. A., Ehrlén, J., Jönsson, M., De Frenne, P., Berg, P., Andersson, J., Weibull, H. & Hylander, N. 2022. ')
$s9_1 = $nr.End
$nr.InsertAfter('Interactive effects of drought and edge exposure on old-growth forest understory species. ')
$e9_1 = $nr.End
$nr.InsertAfter('Landscape Ecology, 37, sid 1839-1853')
$ir = $d.Range($s9_1 - 1, $e9_1 - 1)
$ir.Font.Italic = $true

# --- new paragraph 11 of 13 (style: Normal) ---
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Style = 'Normal'
$nr = $newPara.Range
$nr.InsertAfter('Rudolphi, J., Jönsson, M. T., & Gustafsson, L., 2014. ')
$s10_1 = $nr.End
$nr.InsertAfter('Biological legacies buffer local species extinction after logging. ')
$e10_1 = $nr.End
$nr.InsertAfter('Journal of Applied Ecology. 51, 53-62.')
$ir = $d.Range($s10_1 - 1, $e10_1 - 1)
$ir.Font.Italic = $true

# --- new paragraph 12 of 13 (style: Normal) ---
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Style = 'Normal'
$nr = $newPara.Range
$nr.InsertAfter('Skogsstyrelsen, 2022. ')
$s11_1 = $nr.End
$nr.InsertAfter('Vägledning för hänsyn till knärot. ')
$e11_1 = $nr.End
$nr.InsertAfter('https://www.skogsstyrelsen.se/lag-och-tillsyn/artskydd/vagledningar-och-kunskapsstod-artskydd/vagledning-for-hansyn-till-knarot/')
$ir = $d.Range($s11_1 - 1, $e11_1 - 1)
$ir.Font.Italic = $true

# --- new paragraph 13 of 13 (style: Normal) ---
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Style = 'Normal'
$nr = $newPara.Range
$nr.InsertAfter('SLU Artdatabanken, 2021. ')
$s12_1 = $nr.End
$nr.InsertAfter('Artfaktablad. Naturvård – artfakta. ')
$e12_1 = $nr.End
$nr.InsertAfter('SLU Artdatabanken, Uppsala ')
$ir = $d.Range($s12_1 - 1, $e12_1 - 1)
$ir.Font.Italic = $true

# ----------------------------------------------------------------------
# 2) Update the header date from 2023-09-13 to 2023-09-15.
# ----------------------------------------------------------------------
$d.Content.Find.Execute("2023-09-13", $true, $false, $false, $false, $false,
                        $true, 1, $false, "2023-09-15", 2) | Out-Null

foreach ($sec in $d.Sections) {
    foreach ($idx in 1, 2, 3) {
        $hf = $sec.Headers.Item($idx)
        if ($hf.Exists) {
            $hf.Range.Find.Execute("2023-09-13", $true, $false, $false, $false, $false,
                                   $true, 1, $false, "2023-09-15", 2) | Out-Null
        }
    }
}

Write-Host "done"